# Adds the 2023-12 fixtures (rows 36-44) to the Gibraltar National League sheet,
# matching the data appended to the source workbook on 19-12-2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Duplicate the formatting of the last existing data row (row 35) onto the new
#    rows so column A keeps the bold/centered/bordered index style and column E
#    keeps the custom date-time display style, exactly like every other data row.
$ws.Range("A35:V35").Copy() | Out-Null
$ws.Range("A36:V44").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# 2) Fill in the values for the new rows, cell by cell (keeps number vs. text types
#    exactly as in the source: odds/serial-date columns are numbers, the
#    dd/mm/yyyy hh:mm odds-timestamp columns stay plain text).

# Row 36 (Indice 35): Manchester 62 4-1 College 1975 FC
$ws.Cells.Item(36, 1).Value = 35
$ws.Cells.Item(36, 2).Value = "gibraltar"
$ws.Cells.Item(36, 3).Value = "national-league"
$ws.Cells.Item(36, 4).Value = "2023-2024"
$ws.Cells.Item(36, 5).Value = 45263.6875
$ws.Cells.Item(36, 6).Value = "Manchester 62"
$ws.Cells.Item(36, 7).Value = 4
$ws.Cells.Item(36, 8).Value = "College 1975 FC"
$ws.Cells.Item(36, 9).Value = 1
$ws.Cells.Item(36, 10).Value = 1.34
$ws.Cells.Item(36, 11).Value = "03/12/2023 12:42"
$ws.Cells.Item(36, 12).Value = 1.31
$ws.Cells.Item(36, 13).Value = "03/12/2023 16:29"
$ws.Cells.Item(36, 14).Value = 5.35
$ws.Cells.Item(36, 15).Value = "03/12/2023 12:42"
$ws.Cells.Item(36, 16).Value = 6
$ws.Cells.Item(36, 17).Value = "03/12/2023 16:29"
$ws.Cells.Item(36, 18).Value = 5.51
$ws.Cells.Item(36, 19).Value = "03/12/2023 12:42"
$ws.Cells.Item(36, 20).Value = 5.6
$ws.Cells.Item(36, 21).Value = "03/12/2023 16:29"
$ws.Cells.Item(36, 22).Value = "https://www.betexplorer.com/football/gibraltar/national-league/manchester-62-college-1975/QBPXJmxl/"

# Row 37 (Indice 36): Europa Point 1-0 St Josephs
$ws.Cells.Item(37, 1).Value = 36
$ws.Cells.Item(37, 2).Value = "gibraltar"
$ws.Cells.Item(37, 3).Value = "national-league"
$ws.Cells.Item(37, 4).Value = "2023-2024"
$ws.Cells.Item(37, 5).Value = 45263.8125
$ws.Cells.Item(37, 6).Value = "Europa Point"
$ws.Cells.Item(37, 7).Value = 1
$ws.Cells.Item(37, 8).Value = "St Josephs"
$ws.Cells.Item(37, 9).Value = 0
$ws.Cells.Item(37, 10).Value = 13.15
$ws.Cells.Item(37, 11).Value = "03/12/2023 13:16"
$ws.Cells.Item(37, 12).Value = 14.28
$ws.Cells.Item(37, 13).Value = "03/12/2023 19:01"
$ws.Cells.Item(37, 14).Value = 12.94
$ws.Cells.Item(37, 15).Value = "03/12/2023 13:16"
$ws.Cells.Item(37, 16).Value = 13.12
$ws.Cells.Item(37, 17).Value = "03/12/2023 19:01"
$ws.Cells.Item(37, 18).Value = 1.06
$ws.Cells.Item(37, 19).Value = "03/12/2023 13:16"
$ws.Cells.Item(37, 20).Value = 1.05
$ws.Cells.Item(37, 21).Value = "03/12/2023 18:37"
$ws.Cells.Item(37, 22).Value = "https://www.betexplorer.com/football/gibraltar/national-league/europa-point-st-josephs/42OyJ7if/"

# Row 38 (Indice 37): Magpies 3-0 Glacis United
$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = "gibraltar"
$ws.Cells.Item(38, 3).Value = "national-league"
$ws.Cells.Item(38, 4).Value = "2023-2024"
$ws.Cells.Item(38, 5).Value = 45268.875
$ws.Cells.Item(38, 6).Value = "Magpies"
$ws.Cells.Item(38, 7).Value = 3
$ws.Cells.Item(38, 8).Value = "Glacis United"
$ws.Cells.Item(38, 9).Value = 0
$ws.Cells.Item(38, 10).Value = 1.44
$ws.Cells.Item(38, 11).Value = "08/12/2023 11:01"
$ws.Cells.Item(38, 12).Value = 1.21
$ws.Cells.Item(38, 13).Value = "08/12/2023 20:59"
$ws.Cells.Item(38, 14).Value = 4.55
$ws.Cells.Item(38, 15).Value = "08/12/2023 11:01"
$ws.Cells.Item(38, 16).Value = 5.9
$ws.Cells.Item(38, 17).Value = "08/12/2023 20:59"
$ws.Cells.Item(38, 18).Value = 4.98
$ws.Cells.Item(38, 19).Value = "08/12/2023 11:01"
$ws.Cells.Item(38, 20).Value = 8.76
$ws.Cells.Item(38, 21).Value = "08/12/2023 20:59"
$ws.Cells.Item(38, 22).Value = "https://www.betexplorer.com/football/gibraltar/national-league/magpies-glacis-united/pEuZHR60/"

# Row 39 (Indice 38): College 1975 FC 1-3 Mons Calpe
$ws.Cells.Item(39, 1).Value = 38
$ws.Cells.Item(39, 2).Value = "gibraltar"
$ws.Cells.Item(39, 3).Value = "national-league"
$ws.Cells.Item(39, 4).Value = "2023-2024"
$ws.Cells.Item(39, 5).Value = 45269.6875
$ws.Cells.Item(39, 6).Value = "College 1975 FC"
$ws.Cells.Item(39, 7).Value = 1
$ws.Cells.Item(39, 8).Value = "Mons Calpe"
$ws.Cells.Item(39, 9).Value = 3
$ws.Cells.Item(39, 10).Value = 4.74
$ws.Cells.Item(39, 11).Value = "09/12/2023 12:02"
$ws.Cells.Item(39, 12).Value = 6.36
$ws.Cells.Item(39, 13).Value = "09/12/2023 16:29"
$ws.Cells.Item(39, 14).Value = 4.47
$ws.Cells.Item(39, 15).Value = "09/12/2023 12:02"
$ws.Cells.Item(39, 16).Value = 5.13
$ws.Cells.Item(39, 17).Value = "09/12/2023 16:29"
$ws.Cells.Item(39, 18).Value = 1.47
$ws.Cells.Item(39, 19).Value = "09/12/2023 12:02"
$ws.Cells.Item(39, 20).Value = 1.33
$ws.Cells.Item(39, 21).Value = "09/12/2023 16:29"
$ws.Cells.Item(39, 22).Value = "https://www.betexplorer.com/football/gibraltar/national-league/college-1975-mons-calpe/O4vwHoM6/"

# Row 40 (Indice 39): Europa FC 1-4 Lincoln Red Imps
$ws.Cells.Item(40, 1).Value = 39
$ws.Cells.Item(40, 2).Value = "gibraltar"
$ws.Cells.Item(40, 3).Value = "national-league"
$ws.Cells.Item(40, 4).Value = "2023-2024"
$ws.Cells.Item(40, 5).Value = 45269.8125
$ws.Cells.Item(40, 6).Value = "Europa FC"
$ws.Cells.Item(40, 7).Value = 1
$ws.Cells.Item(40, 8).Value = "Lincoln Red Imps"
$ws.Cells.Item(40, 9).Value = 4
$ws.Cells.Item(40, 10).Value = 19.29
$ws.Cells.Item(40, 11).Value = "09/12/2023 12:04"
$ws.Cells.Item(40, 12).Value = 23.01
$ws.Cells.Item(40, 13).Value = "09/12/2023 19:07"
$ws.Cells.Item(40, 14).Value = 18.58
$ws.Cells.Item(40, 15).Value = "09/12/2023 12:04"
$ws.Cells.Item(40, 16).Value = 20.25
$ws.Cells.Item(40, 17).Value = "09/12/2023 19:07"
$ws.Cells.Item(40, 18).Value = 1.02
$ws.Cells.Item(40, 19).Value = "09/12/2023 12:04"
$ws.Cells.Item(40, 20).Value = 1.02
$ws.Cells.Item(40, 21).Value = "09/12/2023 17:34"
$ws.Cells.Item(40, 22).Value = "https://www.betexplorer.com/football/gibraltar/national-league/europa-fc-lincoln-red-imps/IRwsG5yD/"

# Row 41 (Indice 40): Lynx 1-1 Europa Point
$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 2).Value = "gibraltar"
$ws.Cells.Item(41, 3).Value = "national-league"
$ws.Cells.Item(41, 4).Value = "2023-2024"
$ws.Cells.Item(41, 5).Value = 45270.8125
$ws.Cells.Item(41, 6).Value = "Lynx"
$ws.Cells.Item(41, 7).Value = 1
$ws.Cells.Item(41, 8).Value = "Europa Point"
$ws.Cells.Item(41, 9).Value = 1
$ws.Cells.Item(41, 10).Value = 1.3
$ws.Cells.Item(41, 11).Value = "10/12/2023 11:46"
$ws.Cells.Item(41, 12).Value = 1.31
$ws.Cells.Item(41, 13).Value = "10/12/2023 18:52"
$ws.Cells.Item(41, 14).Value = 5.51
$ws.Cells.Item(41, 15).Value = "10/12/2023 11:46"
$ws.Cells.Item(41, 16).Value = 5.61
$ws.Cells.Item(41, 17).Value = "10/12/2023 19:23"
$ws.Cells.Item(41, 18).Value = 6.27
$ws.Cells.Item(41, 19).Value = "10/12/2023 11:46"
$ws.Cells.Item(41, 20).Value = 6.16
$ws.Cells.Item(41, 21).Value = "10/12/2023 19:23"
$ws.Cells.Item(41, 22).Value = "https://www.betexplorer.com/football/gibraltar/national-league/lynx-europa-point/W8pgNMrt/"

# Row 42 (Indice 41): Manchester 62 5-1 Lynx
$ws.Cells.Item(42, 1).Value = 41
$ws.Cells.Item(42, 2).Value = "gibraltar"
$ws.Cells.Item(42, 3).Value = "national-league"
$ws.Cells.Item(42, 4).Value = "2023-2024"
$ws.Cells.Item(42, 5).Value = 45275.875
$ws.Cells.Item(42, 6).Value = "Manchester 62"
$ws.Cells.Item(42, 7).Value = 5
$ws.Cells.Item(42, 8).Value = "Lynx"
$ws.Cells.Item(42, 9).Value = 1
$ws.Cells.Item(42, 10).Value = 3.93
$ws.Cells.Item(42, 11).Value = "15/12/2023 11:17"
$ws.Cells.Item(42, 12).Value = 2.53
$ws.Cells.Item(42, 13).Value = "15/12/2023 20:59"
$ws.Cells.Item(42, 14).Value = 4.29
$ws.Cells.Item(42, 15).Value = "15/12/2023 11:17"
$ws.Cells.Item(42, 16).Value = 3.7
$ws.Cells.Item(42, 17).Value = "15/12/2023 20:58"
$ws.Cells.Item(42, 18).Value = 1.6
$ws.Cells.Item(42, 19).Value = "15/12/2023 11:17"
$ws.Cells.Item(42, 20).Value = 2.28
$ws.Cells.Item(42, 21).Value = "15/12/2023 20:59"
$ws.Cells.Item(42, 22).Value = "https://www.betexplorer.com/football/gibraltar/national-league/manchester-62-lynx/AaqcMtbn/"

# Row 43 (Indice 42): Europa Point 0-1 Magpies
$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = "gibraltar"
$ws.Cells.Item(43, 3).Value = "national-league"
$ws.Cells.Item(43, 4).Value = "2023-2024"
$ws.Cells.Item(43, 5).Value = 45276.6875
$ws.Cells.Item(43, 6).Value = "Europa Point"
$ws.Cells.Item(43, 7).Value = 0
$ws.Cells.Item(43, 8).Value = "Magpies"
$ws.Cells.Item(43, 9).Value = 1
$ws.Cells.Item(43, 10).Value = 6.71
$ws.Cells.Item(43, 11).Value = "16/12/2023 13:41"
$ws.Cells.Item(43, 12).Value = 8.380000000000001
$ws.Cells.Item(43, 13).Value = "16/12/2023 16:07"
$ws.Cells.Item(43, 14).Value = 5.81
$ws.Cells.Item(43, 15).Value = "16/12/2023 13:41"
$ws.Cells.Item(43, 16).Value = 7.56
$ws.Cells.Item(43, 17).Value = "16/12/2023 16:07"
$ws.Cells.Item(43, 18).Value = 1.25
$ws.Cells.Item(43, 19).Value = "16/12/2023 13:41"
$ws.Cells.Item(43, 20).Value = 1.17
$ws.Cells.Item(43, 21).Value = "16/12/2023 16:07"
$ws.Cells.Item(43, 22).Value = "https://www.betexplorer.com/football/gibraltar/national-league/europa-point-magpies/vZm1L0Dh/"

# Row 44 (Indice 43): Lions Gibraltar 0-5 Europa FC
$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = "gibraltar"
$ws.Cells.Item(44, 3).Value = "national-league"
$ws.Cells.Item(44, 4).Value = "2023-2024"
$ws.Cells.Item(44, 5).Value = 45276.8125
$ws.Cells.Item(44, 6).Value = "Lions Gibraltar"
$ws.Cells.Item(44, 7).Value = 0
$ws.Cells.Item(44, 8).Value = "Europa FC"
$ws.Cells.Item(44, 9).Value = 5
$ws.Cells.Item(44, 10).Value = 3.73
$ws.Cells.Item(44, 11).Value = "16/12/2023 13:41"
$ws.Cells.Item(44, 12).Value = 6.69
$ws.Cells.Item(44, 13).Value = "16/12/2023 19:26"
$ws.Cells.Item(44, 14).Value = 4.24
$ws.Cells.Item(44, 15).Value = "16/12/2023 13:41"
$ws.Cells.Item(44, 16).Value = 4.91
$ws.Cells.Item(44, 17).Value = "16/12/2023 19:26"
$ws.Cells.Item(44, 18).Value = 1.64
$ws.Cells.Item(44, 19).Value = "16/12/2023 13:41"
$ws.Cells.Item(44, 20).Value = 1.33
$ws.Cells.Item(44, 21).Value = "16/12/2023 19:21"
$ws.Cells.Item(44, 22).Value = "https://www.betexplorer.com/football/gibraltar/national-league/lions-gibraltar-europa-fc/IPn5KKSb/"

Write-Host "Inserted rows 36-44 (Indice 35-43)"
